# Separate key viewer and limiter
#
# The "KeyLimiter" sheet used to contain both the key-restriction rows and
# the key-viewer-display rows. This change splits the viewer-display rows
# out into their own new "KeyViewer" sheet (inserted right after
# "KeyLimiter", before "Miscellaneous"), and retitles/retranslates a couple
# of strings along the way:
#   - the old "SHOW_KEY_VIEWER" row becomes the new sheet's "NAME" row
#     ("Key Viewer" / "키뷰어" / "Teclas en pantella")
#   - the old "Show key viewer for registered keys" / "Mostrar teclas en
#     pantalla" row becomes the new sheet's "DESCRIPTION" row ("Shows a key
#     viewer for registered keys" / "등록된 키들의 키뷰어 보이기" /
#     "Muestra teclas en pantalla")
#
# Sheets after "KeyLimiter" are recreated (in the same relative order) so
# that the workbook's internal sheet bookkeeping renumbers cleanly around
# the newly-inserted sheet.

function Get-SheetSnapshot($sheet) {
    $used = $sheet.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $snapshot = New-Object System.Collections.ArrayList
    for ($r = 1; $r -le $rowCount; $r++) {
        $rowData = New-Object System.Collections.ArrayList
        for ($c = 1; $c -le $colCount; $c++) {
            $null = $rowData.Add($sheet.Cells.Item($r, $c).Value2)
        }
        $null = $snapshot.Add($rowData)
    }
    return $snapshot
}

function Write-SheetSnapshot($sheet, $snapshot) {
    for ($r = 0; $r -lt $snapshot.Count; $r++) {
        $rowData = $snapshot[$r]
        for ($c = 0; $c -lt $rowData.Count; $c++) {
            $value = $rowData[$c]
            if ($null -ne $value) {
                $sheet.Cells.Item($r + 1, $c + 1).Value = $value
            }
        }
    }
}

function Get-RowValues($sheet, $row) {
    return @(
        $sheet.Cells.Item($row,1).Value2,
        $sheet.Cells.Item($row,2).Value2,
        $sheet.Cells.Item($row,3).Value2,
        $sheet.Cells.Item($row,4).Value2
    )
}

function Set-RowValues($sheet, $row, $values) {
    $sheet.Cells.Item($row,1).Value = $values[0]
    $sheet.Cells.Item($row,2).Value = $values[1]
    $sheet.Cells.Item($row,3).Value = $values[2]
    $sheet.Cells.Item($row,4).Value = $values[3]
}

$wb = $excel.ActiveWorkbook
$keyLimiter = $wb.Worksheets.Item("KeyLimiter")

# --- 1. Pull out the rows that are moving from KeyLimiter to the new
#        KeyViewer sheet, before anything gets deleted/renumbered. ---
$registeredKeysRow = Get-RowValues $keyLimiter 4
$doneRow = Get-RowValues $keyLimiter 5
$pressKeyRegisterRow = Get-RowValues $keyLimiter 6
$changeKeysRow = Get-RowValues $keyLimiter 7
$viewerOnlyGameplayRow = Get-RowValues $keyLimiter 9
$animateKeysRow = Get-RowValues $keyLimiter 10
$keyViewerSizeRow = Get-RowValues $keyLimiter 11
$keyViewerXPosRow = Get-RowValues $keyLimiter 12
$keyViewerYPosRow = Get-RowValues $keyLimiter 13
$pressedOutlineColorRow = Get-RowValues $keyLimiter 14
$releasedOutlineColorRow = Get-RowValues $keyLimiter 15
$pressedBackgroundColorRow = Get-RowValues $keyLimiter 16
$releasedBackgroundColorRow = Get-RowValues $keyLimiter 17
$pressedTextColorRow = Get-RowValues $keyLimiter 18
$releasedTextColorRow = Get-RowValues $keyLimiter 19

# --- 2. Snapshot the sheets that come after KeyLimiter, then delete them,
#        so they can be recreated (in order) after the new KeyViewer
#        sheet. This keeps the workbook's sheet numbering contiguous. ---
$trailingNames = @("Miscellaneous", "PlanetColor", "PlanetOpacity", "RestrictJudgments")
$trailingSnapshots = @{}
foreach ($name in $trailingNames) {
    $sheet = $wb.Worksheets.Item($name)
    $trailingSnapshots[$name] = Get-SheetSnapshot $sheet
}
foreach ($name in $trailingNames) {
    $wb.Worksheets.Item($name).Delete()
}

# --- 3. Trim KeyLimiter down to just rows 1-7 (header + KEY/NAME/
#        DESCRIPTION/REGISTERED_KEYS/DONE/PRESS_KEY_REGISTER/CHANGE_KEYS). ---
$keyLimiter.Rows("8:19").Delete()

# --- 4. Insert the new KeyViewer sheet right after KeyLimiter. ---
$keyViewer = $wb.Worksheets.Add($null, $keyLimiter)
$keyViewer.Name = "KeyViewer"

Set-RowValues $keyViewer 1 @("KEY", "ENGLISH", "KOREAN", "SPANISH")
Set-RowValues $keyViewer 2 @("NAME", "Key Viewer", "키뷰어", "Teclas en pantella")
Set-RowValues $keyViewer 3 @("DESCRIPTION", "Shows a key viewer for registered keys", "등록된 키들의 키뷰어 보이기", "Muestra teclas en pantalla")

Set-RowValues $keyViewer 4 $registeredKeysRow
Set-RowValues $keyViewer 5 $doneRow
Set-RowValues $keyViewer 6 $pressKeyRegisterRow
Set-RowValues $keyViewer 7 $changeKeysRow
Set-RowValues $keyViewer 8 $viewerOnlyGameplayRow
Set-RowValues $keyViewer 9 $animateKeysRow
Set-RowValues $keyViewer 10 $keyViewerSizeRow
Set-RowValues $keyViewer 11 $keyViewerXPosRow
Set-RowValues $keyViewer 12 $keyViewerYPosRow
Set-RowValues $keyViewer 13 $pressedOutlineColorRow
Set-RowValues $keyViewer 14 $releasedOutlineColorRow
Set-RowValues $keyViewer 15 $pressedBackgroundColorRow
Set-RowValues $keyViewer 16 $releasedBackgroundColorRow
Set-RowValues $keyViewer 17 $pressedTextColorRow
Set-RowValues $keyViewer 18 $releasedTextColorRow

# --- 5. Recreate the trailing sheets (Miscellaneous, PlanetColor,
#        PlanetOpacity, RestrictJudgments) after KeyViewer, in order,
#        restoring their original content. ---
$previous = $keyViewer
foreach ($name in $trailingNames) {
    $newSheet = $wb.Worksheets.Add($null, $previous)
    $newSheet.Name = $name
    Write-SheetSnapshot $newSheet $trailingSnapshots[$name]
    $previous = $newSheet
}
